$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1 = "Save", matching style of existing header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add new data values in column H
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
